# IPI.xlsx upload edit:
#  - fix the spelling of three province names in column A (shared across
#    the 5 year-rows each province occupies)
#  - leave the grid scrolled/selected where the author had it when they
#    stopped editing (B448)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "L' Aquila" -> "L'Aquila"  (rows 7-11)
$ws.Range("A7:A11").Value = "L'Aquila"

# "Reggio nell' Emilia" -> "Reggio nell'Emilia"  (rows 117-121)
$ws.Range("A117:A121").Value = "Reggio nell'Emilia"

# "Massa Carrara" -> "Massa-Carrara"  (rows 447-451)
$ws.Range("A447:A451").Value = "Massa-Carrara"

# Restore the view to where the author left it (scrolled to row 439,
# active cell B448) instead of the stale B536 selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 439
$win.ScrollColumn = 1
$ws.Range("B448").Select()
